$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the data set. It belongs
# right above the current row 34, so insert a blank row there first,
# pushing row 34 (and everything below it) down by one.
$ws.Rows("34:34").Insert()

# Populate the newly inserted row 34 with the new record's data.
$ws.Cells.Item(34, 1).Value = 2
$ws.Cells.Item(34, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(34, 3).Value = "Coquimbo"
$ws.Cells.Item(34, 4).Value = 45162
$ws.Cells.Item(34, 5).Value = 4
$ws.Cells.Item(34, 6).Value = 100112026
$ws.Cells.Item(34, 7).Value = "Haba"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 1300
$ws.Cells.Item(34, 11).Value = 7000
$ws.Cells.Item(34, 12).Value = 8000
$ws.Cells.Item(34, 13).Value = 7500
$ws.Cells.Item(34, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(34, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(34, 16).Value = 300
$ws.Cells.Item(34, 17).Value = 25
$ws.Cells.Item(34, 18).Value = "Hortaliza"

# Match the date cell's number format style used by the rest of column D.
$ws.Cells.Item(34, 4).NumberFormat = $ws.Cells.Item(35, 4).NumberFormat
